$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# NumberFormat is forced to Text ("@") before assignment so that values such as
# "0.9998", "0.07810", or "1.848.22" are preserved exactly as strings instead of
# being auto-coerced into numbers/dates by Excel (matching the original inline-string cells).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.199.50'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.01%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.848.22'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.37%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '246.08'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +1.94%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.81%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9998'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07715'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.55%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3065'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.07%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.59%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07810'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.07%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '93.28'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.48%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.848.31'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.44%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6873'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.22%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.587'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.06%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008324'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.50%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.194.85'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -2.44%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '2.093.14'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.74'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9998'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.517'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9997'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1508'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +0.87%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '159.25'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.852'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -0.91%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.541'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.231'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.179'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.63%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.199'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05118'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.79%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7884'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.91%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.895'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +2.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.149'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.66%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.49%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.321.70'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +7.99%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01868'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.711'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9631'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +7.06%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +9.96%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '106.94'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.43%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '9.706'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.01%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.992.88'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.60%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '64.56'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.37%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.982'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.69%  '
